$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend a "(CODEBLOCK)" marker line to the two code-sample cells (B2, B3).
$jsCode = @'
(CODEBLOCK)
exports.handler = async (event, context, callback) => {
    const hasError = event['queryStringParameters']['myErrorParam'];
    if (hasError === 'yes') {
        callback(new Error('My error message'));
    }
    else {
        const response = {
            statusCode: 200,
            headers: {"Access-Control-Allow-Origin": "*"},
            body: JSON.stringify({success: true}),
            isBase64Encoded: false
        };
        callback(null, response);
}
'@

$batchCode = @'
(CODEBLOCK)
:: This batch file redeploys an existing lambda function
:: Usage: deploy-existing.bat getAllProducts
echo off
echo WARNING: this will delete any index.js or index.zip you have in the current directory!
pause
set /p toDeploy=Enter lambda name (without the .js): 
powershell -Command "(gc %toDeploy%.js) -replace './helpers', '/opt/lambdas/helpers' | Out-File -encoding ASCII index.js"
powershell "Compress-Archive index.js index.zip"
aws lambda update-function-code --function-name %toDeploy% --zip-file fileb://index.zip
del index.js
del index.zip
'@

$ws.Range("B2").Value = $jsCode
$ws.Range("B3").Value = $batchCode

# Row heights grow to fit the extra "(CODEBLOCK)" line now in each cell.
$ws.Rows.Item(2).RowHeight = 217.5
$ws.Rows.Item(3).RowHeight = 188.5

# Move the selection (and scroll position) from A2 to B2.
$ws.Range("B2").Select()
